$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.273.59'
$ws.Range('D3').Value = '3.747.35'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  -0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '593.01'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '165.29'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.32%  '
$ws.Range('D7').Value = '3.747.44'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.07%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.517'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('E10').Value = '  -0.76%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '6.33'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -2.05%  '
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('E13').Value = '  -2.09%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '36.01'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '4.378.71'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '3.734.89'
$ws.Range('E16').Value = '  -0.51%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '18.38'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').Value = '67.268.62'
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('E21').Value = '  -7.45%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '454.28'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.38%  '
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  +5.76%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '83.06'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.79%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.13'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -2.69%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '11.84'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.25%  '
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('E33').Value = '  -0.35%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '9.15'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.63%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '3.700.90'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  -1.17%  '
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -1.17%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.994'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +0.02%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '45.16'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('E45').Value = '  -2.28%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '46.98'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +2.12%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '148.56'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.24%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '8.31'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.96%  '
$ws.Range('E49').Value = '  -4.86%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '388.37'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.41%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '25.94'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.22%  '
